# Auto-generated edit script for cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '29.028.50'
$ws.Range("E2").Value = '  -0.05%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.828.37'
$ws.Range("E3").Value = '  -0.14%  '

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9972'
$ws.Range("E4").Value = '  -0.18%  '

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.20'
$ws.Range("E5").Value = '  +1.13%  '

# Row 6: XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6316'
$ws.Range("E6").Value = '  +0.68%  '

# Row 7: USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").Value = '  -0.03%  '

# Row 8: Dogecoin
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07516'
$ws.Range("E8").Value = '  -1.22%  '

# Row 9: Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2941'
$ws.Range("E9").Value = '  +0.71%  '

# Row 10: Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.09'
$ws.Range("E10").Value = '  +1.23%  '

# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07700'
$ws.Range("E11").Value = '  +0.71%  '

# Row 12: WrappedEther
$ws.Range("D12").Value = '1.831.85'
$ws.Range("E12").Value = '  -0.06%  '

# Row 13: Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.993'
$ws.Range("E13").Value = '  +0.67%  '

# Row 14: Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6690'
$ws.Range("E14").Value = '  +0.50%  '

# Row 15: Litecoin
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.11'
$ws.Range("E15").Value = '  +0.83%  '

# Row 16: ShibaInu
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009614'
$ws.Range("E16").Value = '  +1.80%  '

# Row 17: Uniswap
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.053'
$ws.Range("E17").Value = '  +0.99%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '29.059.49'
$ws.Range("E18").Value = '  +0.35%  '

# Row 19: Avalanche
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.58'
$ws.Range("E19").Value = '  +2.08%  '

# Row 20: BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '226.35'
$ws.Range("E20").Value = '  +0.55%  '

# Row 21: Dai
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9989'
$ws.Range("E21").Value = '  -0.05%  '

# Row 22: Chainlink
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.150'
$ws.Range("E22").Value = '  -1.00%  '

# Row 23: BinanceUSD
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9990'
$ws.Range("E23").Value = '  -0.16%  '

# Row 24: Monero
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.58'
$ws.Range("E24").Value = '  -0.28%  '

# Row 25: Stellar
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1427'
$ws.Range("E25").Value = '  +4.56%  '

# Row 26: Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.515'
$ws.Range("E26").Value = '  +1.12%  '

# Row 27: EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.91'
$ws.Range("E27").Value = '  +0.41%  '

# Row 28: PancakeSwap
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.501'
$ws.Range("E28").Value = '  +0.45%  '

# Row 29: Filecoin
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.143'
$ws.Range("E29").Value = '  +2.09%  '

# Row 30: InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.063'
$ws.Range("E30").Value = '  +0.65%  '

# Row 31: Hedera
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05485'
$ws.Range("E31").Value = '  +5.49%  '

# Row 32: Toncoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.202'
$ws.Range("E32").Value = '  +0.14%  '

# Row 33: LidoDAOToken
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.857'
$ws.Range("E33").Value = '  +0.38%  '

# Row 34: ImmutableX
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7447'
$ws.Range("E34").Value = '  +2.09%  '

# Row 35: ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.137'
$ws.Range("E35").Value = '  -1.40%  '

# Row 36: HuobiToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.648'
$ws.Range("E36").Value = '  +1.53%  '

# Row 37: Maker
$ws.Range("D37").Value = '1.242.68'
$ws.Range("E37").Value = '  -2.49%  '

# Row 38: MXToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.751'
$ws.Range("E38").Value = '  -0.29%  '

# Row 39: VeChain
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01780'
$ws.Range("E39").Value = '  -0.34%  '

# Row 40: FraxShare
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.656'
$ws.Range("E40").Value = '  +2.38%  '

# Row 41: TrustWalletToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9020'
$ws.Range("E41").Value = '  +1.36%  '

# Row 42: PaxDollar
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9997'
$ws.Range("E42").Value = '  -0.09%  '

# Row 43: Quant
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.40'
$ws.Range("E43").Value = '  -0.04%  '

# Row 44: BabyDogeCoin
$ws.Range("B44").Value = 'BabyDogeCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.00000000127'
$ws.Range("E44").Value = '  +5.62%  '

# Row 45: RocketPoolETH
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = '1.979.96'
$ws.Range("E45").Value = '  +0.24%  '

# Row 46: Aave
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.27'
$ws.Range("E46").Value = '  +2.23%  '

# Row 47: Mantle
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5089'
$ws.Range("E47").Value = '  -0.35%  '

# Row 48: TheSandbox
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4058'
$ws.Range("E48").Value = '  +1.92%  '

# Row 49: EnergySwap
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.966'
$ws.Range("E49").Value = '  +1.33%  '

# Row 50: RenderToken
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.655'
$ws.Range("E50").Value = '  +1.07%  '

# Row 51: Cronos
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05785'
$ws.Range("E51").Value = '  +0.72%  '

Write-Output "Update complete"